$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.688.07'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.64%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.512.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.23%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.22'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.33%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.34%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.212'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.25%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.655'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.83%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.33'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.89%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.04%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.44%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.072.22'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.17%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '607.44'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.16%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.744.60'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.57%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.00'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.74%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.58'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.501.77'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.26%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.986'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.55%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.16'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.61%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '106.27'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +12.54%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.08'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.00%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.39%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.01'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.56%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.09%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.79'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.43%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.94%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.54'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.00%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.75%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.22%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.35'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.29%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.14'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.84%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.06%  '

$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.650.40'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.85%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.67'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.76%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.395'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.04%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '509.21'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.18%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.74'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0780'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.80%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.138'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.27%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.41%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.91'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.54%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.21%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.95%  '

$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.75'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.56%  '

$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.94'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.92%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.36'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.51%  '
